$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Jornadas" — Pareja Local / Pareja Visitante columns get the real
# pair names instead of the generic "Pareja N" placeholders.
# ---------------------------------------------------------------------------
$jornadas = $wb.Worksheets.Item("Jornadas")
$jornadas.Activate()

$jornadas.Range("B2").Value = "Jason/Jorge"
$jornadas.Range("C2").Value = "Alex/Keneth"
$jornadas.Range("B3").Value = "Teto/Pedro"
$jornadas.Range("C3").Value = "Keko/Memo"
$jornadas.Range("B4").Value = "Memin/Juan"
$jornadas.Range("C4").Value = "Jason/Jorge"
$jornadas.Range("B5").Value = "Keko/Memo"
$jornadas.Range("C5").Value = "Alex/Keneth"
$jornadas.Range("B6").Value = "Teto/Pedro"
$jornadas.Range("C6").Value = "Keko/Memo"
$jornadas.Range("B7").Value = "Memin/Juan"
$jornadas.Range("C7").Value = "Alex/Keneth"

# Column autofit-style widths recorded in the saved file.
$jornadas.Columns.Item(1).ColumnWidth = 7
$jornadas.Columns.Item(3).ColumnWidth = 14.166666666666666

# New zoom level for this sheet's window.
$excel.ActiveWindow.Zoom = 170

$jornadas.Range("B9").Select()

# ---------------------------------------------------------------------------
# Sheet "Resultados" — same pair-name substitution plus updated game scores.
# ---------------------------------------------------------------------------
$resultados = $wb.Worksheets.Item("Resultados")
$resultados.Activate()

$resultados.Range("B2").Value = "Jason/Jorge"
$resultados.Range("C2").Value = 103
$resultados.Range("D2").Value = "Teto/Pedro"
$resultados.Range("E2").Value = 49

$resultados.Range("B3").Value = "Alex/Keneth"
$resultados.Range("C3").Value = 110
$resultados.Range("D3").Value = "Keko/Memo"
$resultados.Range("E3").Value = 54

$resultados.Range("B4").Value = "Teto/Pedro"
$resultados.Range("C4").Value = 80
$resultados.Range("D4").Value = "Jason/Jorge"
$resultados.Range("E4").Value = 122

$resultados.Range("B5").Value = "Keko/Memo"
$resultados.Range("C5").Value = 120
$resultados.Range("D5").Value = "Alex/Keneth"
$resultados.Range("E5").Value = 76

$resultados.Range("B6").Value = "Jason/Jorge"
$resultados.Range("C6").Value = 110
$resultados.Range("D6").Value = "Keko/Memo"
$resultados.Range("E6").Value = 80

$resultados.Range("B7").Value = "Teto/Pedro"
$resultados.Range("C7").Value = 66
$resultados.Range("D7").Value = "Alex/Keneth"
$resultados.Range("E7").Value = 104

$resultados.Range("E12").Select()

# ---------------------------------------------------------------------------
# Sheet "Clasificación" — only the stored selection rectangle changes.
# ---------------------------------------------------------------------------
$clasificacion = $wb.Worksheets.Item("Clasificación")
$clasificacion.Activate()
$clasificacion.Range("A2:A6").Select()

# Re-activate "Resultados" last so it stays the selected/visible tab,
# matching the workbook's saved active-tab state.
$resultados.Activate()
